$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / workbook view (best effort; cosmetic window placement) ---
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 1440
$win.Width = 37300
$win.Height = 17600

# --- Notes at top of sheet ---
$ws.Range("A1").Value = "This outlines which RASs should be included in the main analysis for a given drug and how they affect the report, based on "
$ws.Range("A2").Value = "(a) the genotype (Gt) & subtype (St) established for the sequence  -- if an unassigned subtype, this counts as ""unknown""."
$ws.Range("A3").Value = "(b) whether the drug is properly characterised in the literature for the genotype & subtypes"

# --- Header row ---
$ws.Range("A5").Value = "Gt known"
$ws.Range("B5").Value = "St known"
$ws.Range("C5").Value = "Gt drug literature"
$ws.Range("D5").Value = "St drug literature"
$ws.Range("E5").Value = "RASs included"
$ws.Range("F5").Value = "Footnote / classification"

# --- Row 7 ---
$ws.Range("A7").Value = "Yes"
$ws.Range("B7").Value = "Yes"
$ws.Range("C7").Value = "Good"
$ws.Range("D7").Value = "Good"
$ws.Range("E7").Value = "Match St or Gt-general"
$ws.Range("F7").Value = "none"

# --- Row 8 ---
$ws.Range("A8").Value = "Yes"
$ws.Range("B8").Value = "No"
$ws.Range("C8").Value = "Good"
$ws.Range("D8").Value = "N/A"
$ws.Range("E8").Value = "Match Gt"
$ws.Range("F8").Value = "Footnote: Since the subtype is unknown, resistant polymorphisms associated with other subtypes of genotype Y were included in the main analysis"

# --- Row 9 ---
$ws.Range("A9").Value = "No"
$ws.Range("B9").Value = "No"
$ws.Range("C9").Value = "N/A"
$ws.Range("D9").Value = "N/A"
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "Classification: Since the genotype is unknown, no resistance analysis was performed"

# --- Row 10 ---
$ws.Range("A10").Value = "Yes"
$ws.Range("B10").Value = "Yes"
$ws.Range("C10").Value = "Good"
$ws.Range("D10").Value = "Poor/None"
$ws.Range("E10").Value = "Match Gt"
$ws.Range("F10").Value = "Footnote: Since the resistance characteristics for drug X in subtype Y are not well understood, resistant polymorphisms associated other subtypes of genotype Z were included in the main analysis"

# --- Row 11 ---
$ws.Range("A11").Value = "Yes"
$ws.Range("B11").Value = "Yes/No"
$ws.Range("C11").Value = "Poor/None"
$ws.Range("D11").Value = "Poor/None"
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = "Classification: ""Since the resistance characteristics for drug X in genotype Y are not well understood, no resistance analysis was performed"

# --- Column widths (character units); runtime quantizes to the nearest
#     pixel internally so these are chosen to land as close as possible
#     to the target stored widths (10, 9, 15.6640625, 15.1640625,
#     21.6640625, 29.33203125). ---
$ws.Columns.Item(1).ColumnWidth = 9.166667
$ws.Columns.Item(2).ColumnWidth = 8.166667
$ws.Columns.Item(3).ColumnWidth = 14.833333
$ws.Columns.Item(4).ColumnWidth = 14.333333
$ws.Columns.Item(5).ColumnWidth = 20.833333
$ws.Columns.Item(6).ColumnWidth = 28.5

# --- View: zoom to 222% and select F11 (last edited cell) ---
$win.Zoom = 222
[void]$ws.Range("F11").Select()
